# Daily scrape update - 2025-08-24 03:27:27 UTC
# Replace the existing data rows with the freshly scraped single opportunity
# and drop the now-stale rows that used to follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 with the new scrape result ---------------------------------
# Column A holds an opportunity id that looks numeric ("1317170"); force it to
# stay text (matches the rest of the sheet, which stores ids as strings) by
# pre-formatting as Text, then strip the leftover number-format so the cell's
# style matches its neighbours again.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1317170"
$ws.Range("A2").ClearFormats() | Out-Null

$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1317170"
$ws.Range("C2").Value = "Guest Relations Officer"
$ws.Range("D2").Value = "Colombo, Sri Lanka"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "39 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Lanka Island Resorts Ltd"

# --- Drop the old rows 3-7 (their data is superseded by the new scrape) ------
$ws.Range("A3:H7").EntireRow.Delete() | Out-Null

# --- Column width changes recorded for this scrape ---------------------------
# Note: this host's ColumnWidth setter round-trips through a px<->chars
# conversion that adds a constant 5/6-character offset vs. the stored OOXML
# <col width>, so the assigned values are pre-compensated (target - 5/6) to
# land on the exact widths (26 / 21 / 15 / 27) recorded in the sheet XML.
$ws.Columns.Item(3).ColumnWidth = 25.166666666666668
$ws.Columns.Item(4).ColumnWidth = 20.166666666666668
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668

Write-Host "UsedRange: $($ws.UsedRange.Address())"
